$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($ws, $ref, $val) {
    $cell = $ws.Range($ref)
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = "Normal"
}

Set-TextCell $ws "D2" "27.450.24"
$ws.Range("E2").Value = "  +0.27%  "

Set-TextCell $ws "D3" "1.737.69"
$ws.Range("E3").Value = "  -0.49%  "

Set-TextCell $ws "D4" "1.002"
$ws.Range("E4").Value = "  -0.46%  "

Set-TextCell $ws "D5" "323.16"
$ws.Range("E5").Value = "  +0.42%  "

Set-TextCell $ws "D6" "1.001"
$ws.Range("E6").Value = "  -0.59%  "

Set-TextCell $ws "D7" "0.4550"
$ws.Range("E7").Value = "  +7.98%  "

$ws.Range("E8").Value = "  -1.70%  "

Set-TextCell $ws "D9" "0.07391"
$ws.Range("E9").Value = "  -1.17%  "

Set-TextCell $ws "D10" "41.24"
$ws.Range("E10").Value = "  -2.96%  "

Set-TextCell $ws "D11" "1.073"
$ws.Range("E11").Value = "  -1.43%  "

Set-TextCell $ws "D12" "1.001"
$ws.Range("E12").Value = "  -0.67%  "

Set-TextCell $ws "D13" "20.37"
$ws.Range("E13").Value = "  -0.85%  "

Set-TextCell $ws "D14" "5.907"
$ws.Range("E14").Value = "  -1.54%  "

Set-TextCell $ws "D15" "7.037"
$ws.Range("E15").Value = "  -2.20%  "

Set-TextCell $ws "D16" "1.735.82"
$ws.Range("E16").Value = "  -1.71%  "

Set-TextCell $ws "D17" "91.15"
$ws.Range("E17").Value = "  +0.56%  "

Set-TextCell $ws "D18" "0.00001052"
$ws.Range("E18").Value = "  -1.44%  "

Set-TextCell $ws "D19" "0.06342"
$ws.Range("E19").Value = "  -0.13%  "

$ws.Range("E20").Value = "  -0.39%  "

Set-TextCell $ws "D21" "16.57"
$ws.Range("E21").Value = "  -2.17%  "

Set-TextCell $ws "D22" "5.716"
$ws.Range("E22").Value = "  -2.48%  "

Set-TextCell $ws "D23" "27.498.64"
$ws.Range("E23").Value = "  +0.24%  "

Set-TextCell $ws "D24" "11.10"
$ws.Range("E24").Value = "  -0.01%  "

Set-TextCell $ws "D25" "2.078"
$ws.Range("E25").Value = "  -0.14%  "

Set-TextCell $ws "D26" "161.84"
$ws.Range("E26").Value = "  +0.25%  "

Set-TextCell $ws "D27" "19.98"
$ws.Range("E27").Value = "  -0.51%  "

Set-TextCell $ws "D28" "1.932.74"
$ws.Range("E28").Value = "  -2.87%  "

Set-TextCell $ws "D29" "2.040"
$ws.Range("E29").Value = "  -3.28%  "

Set-TextCell $ws "D30" "124.60"
$ws.Range("E30").Value = "  +0.71%  "

Set-TextCell $ws "D31" "1.042"
$ws.Range("E31").Value = "  -5.38%  "

Set-TextCell $ws "D32" "0.09081"
$ws.Range("E32").Value = "  +2.58%  "

Set-TextCell $ws "D33" "3.646"
$ws.Range("E33").Value = "  -0.45%  "

Set-TextCell $ws "D34" "5.382"
$ws.Range("E34").Value = "  -2.31%  "

Set-TextCell $ws "D35" "0.02266"
$ws.Range("E35").Value = "  -0.08%  "

Set-TextCell $ws "D36" "11.58"
$ws.Range("E36").Value = "  -5.06%  "

Set-TextCell $ws "D37" "0.05953"
$ws.Range("E37").Value = "  -0.70%  "

Set-TextCell $ws "D38" "0.2056"
$ws.Range("E38").Value = "  -1.39%  "

Set-TextCell $ws "D39" "0.6219"
$ws.Range("E39").Value = "  -1.10%  "

Set-TextCell $ws "D40" "4.872"
$ws.Range("E40").Value = "  -0.64%  "

Set-TextCell $ws "D41" "1.190"
$ws.Range("E41").Value = "  +0.94%  "

Set-TextCell $ws "D42" "1.376"
$ws.Range("E42").Value = "  -1.18%  "

Set-TextCell $ws "D43" "7.678"
$ws.Range("E43").Value = "  -1.97%  "

Set-TextCell $ws "D44" "13.09"
$ws.Range("E44").Value = "  -1.62%  "

Set-TextCell $ws "D45" "3.701"
$ws.Range("E45").Value = "  +0.35%  "

Set-TextCell $ws "D46" "0.5785"
$ws.Range("E46").Value = "  -1.02%  "

Set-TextCell $ws "D47" "121.84"
$ws.Range("E47").Value = "  -0.70%  "

$ws.Range("E48").Value = "  -2.17%  "

$ws.Range("E49").Value = "  +0.46%  "

Set-TextCell $ws "D50" "1.111"
$ws.Range("E50").Value = "  -2.90%  "

Set-TextCell $ws "D51" "71.07"
$ws.Range("E51").Value = "  -2.55%  "
